$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": species counts dropped to 0, percentages cleared ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("B3").Value = 0
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- Sheet "Species qualification": Range Analysis species count dropped to 0 ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# --- Sheet "High Priority break-up": only the IUCN breakdown remains ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Range("A2").Value = "IUCN"
$wsBreak.Range("B2").Value = 18
$wsBreak.Range("C2").Value = 100
$wsBreak.Range("D2").Value = 18
$wsBreak.Range("E2").Value = 100
$wsBreak.Rows("3:3").Delete()
